$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.505.55"
$ws.Cells.Item(2, 5).Value = "  -0.78%  "

$ws.Cells.Item(3, 4).Value = "3.072.70"
$ws.Cells.Item(3, 5).Value = "  -2.20%  "

$ws.Cells.Item(4, 5).Value = "  +0.06%  "

$ws.Cells.Item(5, 4).Value = "'588.81"
$ws.Cells.Item(5, 5).Value = "  -0.34%  "

$ws.Cells.Item(6, 4).Value = "'155.16"
$ws.Cells.Item(6, 5).Value = "  +6.10%  "

$ws.Cells.Item(7, 5).Value = "  -0.05%  "

$ws.Cells.Item(8, 4).Value = "'0.540"
$ws.Cells.Item(8, 5).Value = "  +1.94%  "

$ws.Cells.Item(9, 4).Value = "3.073.39"
$ws.Cells.Item(9, 5).Value = "  -1.86%  "

$ws.Cells.Item(10, 4).Value = "'0.155"
$ws.Cells.Item(10, 5).Value = "  -4.13%  "

$ws.Cells.Item(11, 4).Value = "'5.81"
$ws.Cells.Item(11, 5).Value = "  -1.97%  "

$ws.Cells.Item(12, 4).Value = "'0.451"
$ws.Cells.Item(12, 5).Value = "  -1.12%  "

$ws.Cells.Item(13, 4).Value = "'37.01"
$ws.Cells.Item(13, 5).Value = "  -0.50%  "

$ws.Cells.Item(14, 4).Value = "'0.0000237"
$ws.Cells.Item(14, 5).Value = "  -4.09%  "

$ws.Cells.Item(15, 4).Value = "3.593.06"
$ws.Cells.Item(15, 5).Value = "  -1.77%  "

$ws.Cells.Item(16, 5).Value = "  -1.77%  "

$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "63.616.12"
$ws.Cells.Item(17, 5).Value = "  -0.27%  "

$ws.Cells.Item(18, 2).Value = "Polkadot"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(18, 4).Value = "'7.13"
$ws.Cells.Item(18, 5).Value = "  -1.75%  "

$ws.Cells.Item(19, 4).Value = "3.081.48"
$ws.Cells.Item(19, 5).Value = "  -1.70%  "

$ws.Cells.Item(20, 4).Value = "'475.68"
$ws.Cells.Item(20, 5).Value = "  +2.00%  "

$ws.Cells.Item(21, 4).Value = "'14.41"
$ws.Cells.Item(21, 5).Value = "  +0.45%  "

$ws.Cells.Item(22, 4).Value = "'0.707"
$ws.Cells.Item(22, 5).Value = "  -3.26%  "

$ws.Cells.Item(23, 4).Value = "'7.50"
$ws.Cells.Item(23, 5).Value = "  -0.94%  "

$ws.Cells.Item(24, 4).Value = "'2.40"
$ws.Cells.Item(24, 5).Value = "  +0.18%  "

$ws.Cells.Item(25, 4).Value = "'80.86"
$ws.Cells.Item(25, 5).Value = "  -0.05%  "

$ws.Cells.Item(26, 4).Value = "'12.76"
$ws.Cells.Item(26, 5).Value = "  -2.95%  "

$ws.Cells.Item(27, 4).Value = "'10.25"
$ws.Cells.Item(27, 5).Value = "  +4.42%  "

$ws.Cells.Item(28, 4).Value = "'1.00"
$ws.Cells.Item(28, 5).Value = "  -0.08%  "

$ws.Cells.Item(29, 4).Value = "'7.52"
$ws.Cells.Item(29, 5).Value = "  +2.36%  "

$ws.Cells.Item(30, 4).Value = "'2.66"
$ws.Cells.Item(30, 5).Value = "  -1.91%  "

$ws.Cells.Item(31, 5).Value = "  -0.19%  "

$ws.Cells.Item(32, 4).Value = "'2.15"
$ws.Cells.Item(32, 5).Value = "  -2.51%  "

$ws.Cells.Item(33, 4).Value = "'0.112"
$ws.Cells.Item(33, 5).Value = "  -1.31%  "

$ws.Cells.Item(34, 4).Value = "'27.07"
$ws.Cells.Item(34, 5).Value = "  -2.06%  "

$ws.Cells.Item(35, 4).Value = "0.0₃0819"
$ws.Cells.Item(35, 5).Value = "  -5.07%  "

$ws.Cells.Item(36, 2).Value = "Mantle"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(36, 4).Value = "'1.05"
$ws.Cells.Item(36, 5).Value = "  -1.24%  "

$ws.Cells.Item(37, 2).Value = "dogwifhat"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(37, 4).Value = "'3.35"
$ws.Cells.Item(37, 5).Value = "  +2.32%  "

$ws.Cells.Item(38, 4).Value = "'6.01"
$ws.Cells.Item(38, 5).Value = "  -2.37%  "

$ws.Cells.Item(39, 4).Value = "'2.20"
$ws.Cells.Item(39, 5).Value = "  -3.12%  "

$ws.Cells.Item(40, 4).Value = "'9.32"
$ws.Cells.Item(40, 5).Value = "  -0.54%  "

$ws.Cells.Item(41, 4).Value = "'50.67"
$ws.Cells.Item(41, 5).Value = "  -1.27%  "

$ws.Cells.Item(42, 4).Value = "'442.38"
$ws.Cells.Item(42, 5).Value = "  -4.24%  "

$ws.Cells.Item(43, 4).Value = "'0.287"
$ws.Cells.Item(43, 5).Value = "  -2.14%  "

$ws.Cells.Item(44, 4).Value = "'40.99"
$ws.Cells.Item(44, 5).Value = "  +1.90%  "

$ws.Cells.Item(45, 5).Value = "  +4.40%  "

$ws.Cells.Item(46, 4).Value = "'0.0359"
$ws.Cells.Item(46, 5).Value = "  -3.50%  "

$ws.Cells.Item(47, 4).Value = "2.800.18"
$ws.Cells.Item(47, 5).Value = "  -3.12%  "

$ws.Cells.Item(48, 4).Value = "'130.61"
$ws.Cells.Item(48, 5).Value = "  -1.72%  "

$ws.Cells.Item(49, 2).Value = "USDe"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(49, 4).Value = "'1.00"
$ws.Cells.Item(49, 5).Value = "  +0.09%  "

$ws.Cells.Item(50, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(50, 4).Value = "'25.13"
$ws.Cells.Item(50, 5).Value = "  +4.67%  "

$ws.Cells.Item(51, 4).Value = "'2.24"
$ws.Cells.Item(51, 5).Value = "  +0.49%  "
